$wb = $excel.ActiveWorkbook

# ALC row 19 (anchor G19=7015)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1345.875
$ws.Range("I19").Value = 1323.8572
$ws.Range("J19").Value = 1500
$ws.Range("K19").Value = 1323.8572
$ws.Range("L19").Value = 1500
$ws.Range("M19").Value = -1148.8572
$ws.Range("N19").Value = -1850

# ALC row 28 (anchor G28=27772)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 894.1579
$ws.Range("I28").Value = 694.375
$ws.Range("K28").Value = 694.375
$ws.Range("M28").Value = -209.375

# ALC row 33 (anchor G33=5512)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 333.92307
$ws.Range("I33").Value = 379.63635
$ws.Range("K33").Value = 379.63635
$ws.Range("M33").Value = -150.63635

# ALC row 55 (anchor G55=5517)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 572.625
$ws.Range("I55").Value = 363.8
$ws.Range("J55").Value = 920.6667
$ws.Range("K55").Value = 363.8
$ws.Range("L55").Value = 920.6667
$ws.Range("M55").Value = -149.8
$ws.Range("N55").Value = -1348.6667

# ALC row 132 (anchor G132=44049)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 5932.9375
$ws.Range("I132").Value = 6088.7334
$ws.Range("K132").Value = 18266.2002
$ws.Range("M132").Value = -15736.2002

# ARM row 45 (anchor G45=27714)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1562.1428
$ws.Range("I45").Value = 1655.8334
$ws.Range("J45").Value = 1000
$ws.Range("K45").Value = 1655.8334
$ws.Range("L45").Value = 1000
$ws.Range("M45").Value = -1278.8334
$ws.Range("N45").Value = -1754

# ARM row 76 (anchor G76=10679)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").ClearContents()

# ARM row 79 (anchor G79=10679)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").ClearContents()

# ARM row 110 (anchor G110=27708)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 2418.4443
$ws.Range("I110").Value = 2418.4443
$ws.Range("K110").Value = 2418.4443
$ws.Range("M110").Value = -373.4443000000001

# ARM row 114 (anchor G114=25968)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H114").Value = 50000
$ws.Range("J114").Value = 50000
$ws.Range("L114").Value = 50000
$ws.Range("N114").Value = -58678

# ARM row 132 (anchor G132=43997)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 5551.3335
$ws.Range("I132").Value = 5551.3335
$ws.Range("K132").Value = 16654.0005
$ws.Range("M132").Value = -14124.0005

# BSM row 22 (anchor G22=5092)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 316.2
$ws.Range("I22").Value = 316.2
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 316.2
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -143.2
$ws.Range("N22").ClearContents()

# BSM row 99 (anchor G99=19943)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2000
$ws.Range("I99").Value = 2000
$ws.Range("K99").Value = 2000
$ws.Range("M99").Value = -502

# BSM row 105 (anchor G105=19947)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3813.0435
$ws.Range("I105").Value = 2906.2307
$ws.Range("J105").Value = 4991.9
$ws.Range("K105").Value = 2906.2307
$ws.Range("L105").Value = 4991.9
$ws.Range("M105").Value = -1159.2307
$ws.Range("N105").Value = -8485.9

# BSM row 107 (anchor G107=27706)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 0
$ws.Range("I107").Value = 0
$ws.Range("K107").Value = 0
$ws.Range("M107").ClearContents()

# BSM row 134 (anchor G134=43998)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3497.25
$ws.Range("I134").Value = 3829.6667
$ws.Range("K134").Value = 11489.0001
$ws.Range("M134").Value = -8954.000100000001

# CRP row 31 (anchor G31=44023)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1278.25
$ws.Range("I31").Value = 1278.25
$ws.Range("K31").Value = 1278.25
$ws.Range("M31").Value = -983.25

# CRP row 34 (anchor G34=44023)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 1278.25
$ws.Range("I34").Value = 1278.25
$ws.Range("K34").Value = 1278.25
$ws.Range("M34").Value = -1076.25

# CRP row 58 (anchor G58=44021)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 5025.875
$ws.Range("I58").Value = 5034.6665
$ws.Range("K58").Value = 5034.6665
$ws.Range("M58").Value = -4831.6665

# CRP row 99 (anchor G99=36198)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 2688.6667
$ws.Range("I99").Value = 2712.25
$ws.Range("K99").Value = 2712.25
$ws.Range("M99").Value = -1214.25

# CRP row 107 (anchor G107=27689)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 1641
$ws.Range("I107").Value = 999.5
$ws.Range("J107").Value = 2282.5
$ws.Range("K107").Value = 999.5
$ws.Range("L107").Value = 2282.5
$ws.Range("M107").Value = 920.5
$ws.Range("N107").Value = -6122.5

# CRP row 126 (anchor G126=36198)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 2688.6667
$ws.Range("I126").Value = 2712.25
$ws.Range("K126").Value = 8136.75
$ws.Range("M126").Value = -5666.75

# CRP row 132 (anchor G132=44019)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2883.3333
$ws.Range("I132").Value = 2883.3333
$ws.Range("K132").Value = 8649.999899999999
$ws.Range("M132").Value = -6119.999899999999

# CRP row 134 (anchor G134=44020)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 2638.4167
$ws.Range("I134").Value = 2638.4167
$ws.Range("K134").Value = 7915.250100000001
$ws.Range("M134").Value = -5380.250100000001

# CRP row 136 (anchor G136=44021)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 5025.875
$ws.Range("I136").Value = 5034.6665
$ws.Range("K136").Value = 15103.9995
$ws.Range("M136").Value = -12553.9995

# CUL row 15 (anchor G15=4661)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H15").Value = 417
$ws.Range("J15").Value = 845
$ws.Range("L15").Value = 2535
$ws.Range("N15").Value = -2815

# CUL row 26 (anchor G26=4746)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 1281.8966
$ws.Range("J26").Value = 397.5
$ws.Range("L26").Value = 1192.5
$ws.Range("N26").Value = -1768.5

# CUL row 50 (anchor G50=4725)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H50").Value = 254
$ws.Range("I50").Value = 254
$ws.Range("K50").Value = 762
$ws.Range("M50").Value = -281

# CUL row 53 (anchor G53=4725)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H53").Value = 254
$ws.Range("I53").Value = 254
$ws.Range("K53").Value = 762
$ws.Range("M53").Value = -281

# CUL row 60 (anchor G60=4750)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H60").Value = 2587.111
$ws.Range("J60").Value = 3430
$ws.Range("L60").Value = 10290
$ws.Range("N60").Value = -10792

# CUL row 62 (anchor G62=12845)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H62").Value = 8548.5
$ws.Range("J62").Value = 6497
$ws.Range("L62").Value = 19491
$ws.Range("N62").Value = -20863

# CUL row 65 (anchor G65=12845)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H65").Value = 8548.5
$ws.Range("J65").Value = 6497
$ws.Range("L65").Value = 58473
$ws.Range("N65").Value = -65337

# CUL row 86 (anchor G86=12892)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 1066.1428
$ws.Range("I86").Value = 330
$ws.Range("K86").Value = 990
$ws.Range("M86").Value = 196

# CUL row 89 (anchor G89=12892)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H89").Value = 1066.1428
$ws.Range("I89").Value = 330
$ws.Range("K89").Value = 2970
$ws.Range("M89").Value = 2958

# CUL row 133 (anchor G133=44073)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H133").Value = 16705.883
$ws.Range("I133").Value = 12857.143
$ws.Range("K133").Value = 38571.429
$ws.Range("M133").Value = -33511.429

# GSM row 44 (anchor G44=4143)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H44").Value = 17499.834
$ws.Range("I44").Value = 17499.834
$ws.Range("J44").Value = 0
$ws.Range("K44").Value = 17499.834
$ws.Range("L44").Value = 0
$ws.Range("M44").Value = -16903.834
$ws.Range("N44").ClearContents()

# GSM row 46 (anchor G46=2078)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 3722
$ws.Range("I46").Value = 3722
$ws.Range("K46").Value = 3722
$ws.Range("M46").Value = -3566

# GSM row 58 (anchor G58=4363)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H58").Value = 37500
$ws.Range("I58").Value = 37500
$ws.Range("K58").Value = 37500
$ws.Range("M58").Value = -37223

# GSM row 113 (anchor G113=27710)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 2997.2856
$ws.Range("I113").Value = 2997.8
$ws.Range("J113").Value = 2996
$ws.Range("K113").Value = 2997.8
$ws.Range("L113").Value = 2996
$ws.Range("M113").Value = -827.8000000000002
$ws.Range("N113").Value = -7336

# GSM row 132 (anchor G132=44008)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2465.6667
$ws.Range("I132").Value = 1958.8
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 5876.4
$ws.Range("L132").Value = 15000
$ws.Range("M132").Value = -3346.4
$ws.Range("N132").Value = -20060

# LTW row 40 (anchor G40=36248)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 1301.3334
$ws.Range("I40").Value = 1323.5714
$ws.Range("K40").Value = 1323.5714
$ws.Range("M40").Value = -1187.5714

# LTW row 61 (anchor G61=27740)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1949
$ws.Range("I61").Value = 1949
$ws.Range("K61").Value = 1949
$ws.Range("M61").Value = -1747

# LTW row 100 (anchor G100=19995)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 400
$ws.Range("I100").Value = 400
$ws.Range("K100").Value = 400
$ws.Range("M100").Value = 141

# LTW row 113 (anchor G113=27740)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 1949
$ws.Range("I113").Value = 1949
$ws.Range("K113").Value = 1949
$ws.Range("M113").Value = 221

# LTW row 122 (anchor G122=36247)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 4886.778
$ws.Range("I122").Value = 4881.8335
$ws.Range("J122").Value = 4896.6665
$ws.Range("K122").Value = 14645.5005
$ws.Range("L122").Value = 14689.9995
$ws.Range("M122").Value = -12195.5005
$ws.Range("N122").Value = -19589.9995

# LTW row 132 (anchor G132=44058)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2894.375
$ws.Range("I132").Value = 2894.375
$ws.Range("K132").Value = 8683.125
$ws.Range("M132").Value = -6153.125

# LTW row 136 (anchor G136=44060)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 1056.5
$ws.Range("I136").Value = 1056.5
$ws.Range("K136").Value = 3169.5
$ws.Range("M136").Value = -619.5

# WVR row 107 (anchor G107=27746)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 4070.6667
$ws.Range("I107").Value = 2266.8
$ws.Range("K107").Value = 6800.400000000001
$ws.Range("M107").Value = -4880.400000000001

# WVR row 113 (anchor G113=27752)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 2433
$ws.Range("I113").Value = 2433
$ws.Range("K113").Value = 7299
$ws.Range("M113").Value = -5129

# WVR row 122 (anchor G122=36208)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3713.6956
$ws.Range("I122").Value = 4089.4707
$ws.Range("K122").Value = 12268.4121
$ws.Range("M122").Value = -9818.4121

# WVR row 132 (anchor G132=44029)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2718.5293
$ws.Range("I132").Value = 3040.4
$ws.Range("J132").Value = 304.5
$ws.Range("K132").Value = 9121.200000000001
$ws.Range("L132").Value = 913.5
$ws.Range("M132").Value = -6591.200000000001
$ws.Range("N132").Value = -5973.5
